$wb = $excel.ActiveWorkbook

# --- studies sheet (sheet1): add "PMID" header in column H ---
$wsStudies = $wb.Worksheets.Item("studies")
$wsStudies.Range("H1").Value = "PMID"
$wsStudies.Range("H1").Style = $wsStudies.Range("G1").Style
$wsStudies.Range("H2").Select() | Out-Null

# --- counts sheet (sheet3): add "notes" header in column F ---
$wsCounts = $wb.Worksheets.Item("counts")
$wsCounts.Range("F1").Value = "notes"
$wsCounts.Columns.Item(3).ColumnWidth = 14.4166666669771
$wsCounts.Range("F2").Select() | Out-Null
